$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Vtn", "Itga5", "ECs", 3, 1, 12.75206033333333, 38.256181, 0.1573122343381959, 0.157312234338196, 3, 1, 30.99161333333333, 92.97484, 0.3599121977633812, 0.3599121977633811, 395.2069230540044, 3556.86230748604, 0.05661859199572815, 0.05661859199572814),
    @("ECs", "Vtn", "Itga5", "FAPs", 3, 1, 12.75206033333333, 38.256181, 0.1573122343381959, 0.157312234338196, 3, 1, 29.913269, 89.739807, 0.3473891556493311, 0.3473891556493311, 381.4558110552296, 3433.102299497067, 0.05464856426005561, 0.05464856426005561),
    @("ECs", "Vtn", "Itga5", "M2", 3, 1, 12.75206033333333, 38.256181, 0.1573122343381959, 0.157312234338196, 3, 1, 19.150218, 57.450654, 0.2223955550134164, 0.2223955550134163, 244.204735332486, 2197.842617992374, 0.03498554166604371, 0.03498554166604371),
    @("ECs", "Vtn", "Itga5", "sCs", 3, 1, 12.75206033333333, 38.256181, 0.1573122343381959, 0.157312234338196, 3, 1, 6.053716000000001, 18.161148, 0.07030309157387134, 0.07030309157387132, 77.19735167286534, 694.7761650557879, 0.0110595364163685, 0.0110595364163685),
    @("FAPs", "Vtn", "Itga5", "ECs", 3, 1, 20.35396833333334, 61.06190500000001, 0.2510910513649196, 0.2510910513649196, 3, 1, 30.99161333333333, 92.97484, 0.3599121977633812, 0.3599121977633811, 630.8023163855779, 5677.220847470201, 0.09037073213546626, 0.09037073213546623),
    @("FAPs", "Vtn", "Itga5", "FAPs", 3, 1, 20.35396833333334, 61.06190500000001, 0.2510910513649196, 0.2510910513649196, 3, 1, 29.913269, 89.739807, 0.3473891556493311, 0.3473891556493311, 608.8537299724818, 5479.683569752336, 0.08722630832476226, 0.08722630832476225),
    @("FAPs", "Vtn", "Itga5", "M2", 3, 1, 20.35396833333334, 61.06190500000001, 0.2510910513649196, 0.2510910513649196, 3, 1, 19.150218, 57.450654, 0.2223955550134164, 0.2223955550134163, 389.7829307484301, 3508.046376735871, 0.05584153372720354, 0.05584153372720353),
    @("FAPs", "Vtn", "Itga5", "sCs", 3, 1, 20.35396833333334, 61.06190500000001, 0.2510910513649196, 0.2510910513649196, 3, 1, 6.053716000000001, 18.161148, 0.07030309157387134, 0.07030309157387132, 123.2171437629934, 1108.95429386694, 0.01765247717748757, 0.01765247717748757),
    @("M2", "Vtn", "Itga5", "ECs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.00005805425597465284, 0.00005805425597465285, 3, 1, 30.99161333333333, 92.97484, 0.3599121977633812, 0.3599121977633811, 0.1458465323466667, 1.31261879112, 0.00002089443485735521, 0.00002089443485735521),
    @("M2", "Vtn", "Itga5", "FAPs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.00005805425597465284, 0.00005805425597465285, 3, 1, 29.913269, 89.739807, 0.3473891556493311, 0.3473891556493311, 0.140771843914, 1.266946595226, 0.00002016741896488479, 0.00002016741896488479),
    @("M2", "Vtn", "Itga5", "M2", 1, 0.3333333333333333, 0.004706, 0.014118, 0.00005805425597465284, 0.00005805425597465285, 3, 1, 19.150218, 57.450654, 0.2223955550134164, 0.2223955550134163, 0.090120925908, 0.811088333172, 0.00001291100847837386, 0.00001291100847837386),
    @("M2", "Vtn", "Itga5", "sCs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.00005805425597465284, 0.00005805425597465285, 3, 1, 6.053716000000001, 18.161148, 0.07030309157387134, 0.07030309157387132, 0.028488787496, 0.256399087464, 0.000004081393674038986, 0.000004081393674038986),
    @("sCs", "Vtn", "Itga5", "ECs", 3, 1, 47.951367, 143.854101, 0.5915386600409097, 0.5915386600409098, 3, 1, 30.99161333333333, 92.97484, 0.3599121977633812, 0.3599121977633811, 1486.09022486876, 13374.81202381884, 0.2129019791973294, 0.2129019791973294),
    @("sCs", "Vtn", "Itga5", "FAPs", 3, 1, 47.951367, 143.854101, 0.5915386600409097, 0.5915386600409098, 3, 1, 29.913269, 89.739807, 0.3473891556493311, 0.3473891556493311, 1434.382139988723, 12909.43925989851, 0.2054941156455484, 0.2054941156455484),
    @("sCs", "Vtn", "Itga5", "M2", 3, 1, 47.951367, 143.854101, 0.5915386600409097, 0.5915386600409098, 3, 1, 19.150218, 57.450654, 0.2223955550134164, 0.2223955550134163, 918.2791314480061, 8264.512183032055, 0.1315555686116908, 0.1315555686116907),
    @("sCs", "Vtn", "Itga5", "sCs", 3, 1, 47.951367, 143.854101, 0.5915386600409097, 0.5915386600409098, 3, 1, 6.053716000000001, 18.161148, 0.07030309157387134, 0.07030309157387132, 290.283957629772, 2612.555618667948, 0.04158699658634122, 0.04158699658634122)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}
